# "paises.xlsx" data refresh: 3-Apr-2020 13:50 -> 14:20 snapshot.
#  - Timestamp caption updated.
#  - Most-affected countries (Alemania, Paises Bajos, Austria, Portugal,
#    Brasil, Suecia, Noruega, Principado de Andorra) get new totals only.
#  - Croacia/Emiratos Arabes Unidos (rows 54-55), Kazajistan/Tunez
#    (rows 76-77) and Puerto Rico/Zambia (rows 141-142) additionally swap
#    list order, so both the country name and the figures move per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1 - timestamp caption
$ws.Range("A1").Value = 'Datos actualizados a 3 de Abril de 2020 a las 14:20'

# Row 7
$ws.Range("B7").Value = 85903
$ws.Range("C7").Value = 1109
$ws.Range("E7").Value = 62341
$ws.Range("G7").Value = 15
$ws.Range("H7").Value = 1122

# Row 15
$ws.Range("B15").Value = 15723
$ws.Range("C15").Value = 1026
$ws.Range("E15").Value = 13986
$ws.Range("F15").Value = 1273
$ws.Range("G15").Value = 148
$ws.Range("H15").Value = 1487

# Row 16
$ws.Range("B16").Value = 11350
$ws.Range("C16").Value = 221
$ws.Range("E16").Value = 9160

# Row 19
$ws.Range("B19").Value = 9886
$ws.Range("C19").Value = 852
$ws.Range("E19").Value = 9572
$ws.Range("F19").Value = 245
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = 246

# Row 20
$ws.Range("B20").Value = 8076
$ws.Range("C20").Value = 32
$ws.Range("E20").Value = 7622

# Row 22
$ws.Range("B22").Value = 6078
$ws.Range("C22").Value = 510
$ws.Range("E22").Value = 5642
$ws.Range("F22").Value = 469
$ws.Range("G22").Value = 25
$ws.Range("H22").Value = 333

# Row 24
$ws.Range("B24").Value = 5296
$ws.Range("C24").Value = 149
$ws.Range("E24").Value = 5210

# Row 54 (was Emiratos Arabes Unidos) -> Croacia, with updated figures
$ws.Range("A54").Value = 'Croacia'
$ws.Range("B54").Value = 1079
$ws.Range("C54").Value = 68
$ws.Range("D54").Value = 92
$ws.Range("E54").Value = 979
$ws.Range("F54").Value = 39
$ws.Range("G54").Value = 1

# Row 55 (was Croacia) -> Emiratos Arabes Unidos, figures unchanged
$ws.Range("A55").Value = 'Emiratos Arabes Unidos'
$ws.Range("B55").Value = 1024
$ws.Range("D55").Value = 96
$ws.Range("E55").Value = 920
$ws.Range("F55").Value = 2
$ws.Range("H55").Value = 8

# Row 76 (was Tunez) -> Kazajistan, with updated figures
$ws.Range("A76").Value = 'Kazajistan'
$ws.Range("B76").Value = 460
$ws.Range("C76").Value = 25
$ws.Range("D76").Value = 29
$ws.Range("E76").Value = 425
$ws.Range("F76").Value = 6
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = 6

# Row 77 (was Kazajistan) -> Tunez, figures unchanged
$ws.Range("A77").Value = 'Tunez'
$ws.Range("B77").Value = 455
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 5
$ws.Range("E77").Value = 436
$ws.Range("F77").Value = 10
$ws.Range("H77").Value = 14

# Row 80 - Principado de Andorra, updated figures (no reorder)
$ws.Range("B80").Value = 439
$ws.Range("C80").Value = 11
$ws.Range("D80").Value = 16
$ws.Range("E80").Value = 407
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 16

# Row 141 (was Zambia) -> Puerto Rico, figures unchanged
$ws.Range("A141").Value = 'Puerto Rico'
$ws.Range("D141").Value = 1
$ws.Range("E141").Value = 36
$ws.Range("H141").Value = 2

# Row 142 (was Puerto Rico) -> Zambia, with updated figures
$ws.Range("A142").Value = 'Zambia'
$ws.Range("D142").Value = 2
$ws.Range("H142").Value = 1
